# Refresh the cryptos list table with the latest scraped price / 1h-volume data,
# keeping the existing text formatting (values are stored as plain text, exactly as
# the original scraper wrote them - e.g. "36.505.25" or "  -0.36%  ").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $value into $cellRef as literal text, even when the text looks like
# a plain number (e.g. "251.82"), so Excel does not silently convert it to a Double
# and mangle it with floating point noise (251.81999999999999) or reformat it.
function Set-TextValue([string]$cellRef, [string]$value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}


$ws.Range("D2").Value = "36.505.25"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "2.096.62"
$ws.Range("E3").Value = "  +9.48%  "

$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue "D5" "251.82"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("E6").Value = "  -6.77%  "

$ws.Range("E7").Value = "  +0.08%  "

Set-TextValue "D8" "47.23"
$ws.Range("E8").Value = "  +5.22%  "

Set-TextValue "D9" "59.61"
$ws.Range("E9").Value = "  +1.99%  "

Set-TextValue "D10" "0.373"
$ws.Range("E10").Value = "  +0.79%  "

Set-TextValue "D11" "0.0741"
$ws.Range("E11").Value = "  -3.05%  "

Set-TextValue "D12" "0.0999"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "2.399.77"
$ws.Range("E14").Value = "  +9.49%  "

$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").Value = "2.095.37"
$ws.Range("E16").Value = "  +9.45%  "

$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("D18").Value = "36.512.38"
$ws.Range("E18").Value = "  -0.22%  "

Set-TextValue "D19" "72.57"
$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("D20").Value = "0.0₃0825"
$ws.Range("E20").Value = "  -4.43%  "

Set-TextValue "D21" "13.10"
$ws.Range("E21").Value = "  -2.30%  "

Set-TextValue "D22" "239.08"
$ws.Range("E22").Value = "  -4.56%  "

$ws.Range("E23").Value = "  -1.70%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  -7.02%  "

Set-TextValue "D26" "169.97"
$ws.Range("E26").Value = "  +0.63%  "

Set-TextValue "D27" "21.26"
$ws.Range("E27").Value = "  +13.48%  "

$ws.Range("E28").Value = "  +3.14%  "

Set-TextValue "D29" "1.98"
$ws.Range("E29").Value = "  -9.89%  "

Set-TextValue "D30" "28.12"
$ws.Range("E30").Value = "  +58.21%  "

Set-TextValue "D31" "0.123"
$ws.Range("E31").Value = "  -5.50%  "

$ws.Range("E32").Value = "  -2.91%  "

Set-TextValue "D33" "0.0608"

Set-TextValue "D34" "0.0922"
$ws.Range("E34").Value = "  +3.36%  "

Set-TextValue "D35" "0.977"
$ws.Range("E35").Value = "  +10.39%  "

Set-TextValue "D36" "2.36"
$ws.Range("E36").Value = "  +16.41%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -2.19%  "

Set-TextValue "D39" "4.06"
$ws.Range("E39").Value = "  -6.67%  "

Set-TextValue "D40" "1.34"
$ws.Range("E40").Value = "  -12.16%  "

Set-TextValue "D43" "97.14"
$ws.Range("E43").Value = "  -8.42%  "

$ws.Range("E44").Value = "  -7.00%  "

Set-TextValue "D45" "15.92"
$ws.Range("E45").Value = "  -9.53%  "

$ws.Range("D46").Value = "1.326.60"
$ws.Range("E46").Value = "  -1.37%  "

Set-TextValue "D47" "0.0840"
$ws.Range("E47").Value = "  +3.27%  "

Set-TextValue "D48" "6.94"
$ws.Range("E48").Value = "  +8.42%  "

$ws.Range("D49").Value = "2.288.84"
$ws.Range("E49").Value = "  +9.52%  "

Set-TextValue "D50" "2.84"
$ws.Range("E50").Value = "  +1.48%  "

Set-TextValue "D51" "2.23"
$ws.Range("E51").Value = "  -6.51%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0221"
$ws.Range("E41").Value = "  -2.60%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "1.16"
$ws.Range("E42").Value = "  +4.58%  "
